$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Advanced settings")

$ws.Range("A7").Value = "Enable grid import fee"
$ws.Range("A8").Value = "Enable grid export fee"
$ws.Range("A9").Value = "Enable taxes surcharges"
$ws.Range("A10").Value = "Enable marketplace monthly fee"
$ws.Range("A11").Value = "Enable assistance monthly fee"
$ws.Range("A12").Value = "Enable service monthly fee"
$ws.Range("A13").Value = "Enable contracted power monthly fee"
$ws.Range("A14").Value = "Enable contracted power cargo monthly fee"
$ws.Range("A15").Value = "Enable energy cargo fee"

$ws.Activate()
$ws.Range("A6").Select()
